# Reset the shop's product data: remove the "prod-001" (iPhone 14 Pro) row
# and the "prod-005" (Dell XPS 13) row from the Products sheet, letting the
# remaining rows (prod-002, prod-003, prod-004) shift up so the table is
# contiguous again (A1:L4).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Products")

# Delete row 2 (prod-001 / iPhone 14 Pro). Remaining rows shift up one.
$ws.Rows.Item(2).EntireRow.Delete()

# After the shift, the former row 6 (prod-005 / Dell XPS 13) is now row 5.
$ws.Rows.Item(5).EntireRow.Delete()
